# run initial lines for week 13
# Appends Week 13 betting lines to the "Sheet1" worksheet (rows 164-179),
# matching the layout of: week, game, total_line, spread_line

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$week13 = @(
    @(13, "GB_DET",  48.5,  2.5),
    @(13, "KC_DAL",  47.5, -5.5),
    @(13, "CIN_BAL", 49.5,  5.5),
    @(13, "CHI_PHI", 46.5,  6.5),
    @(13, "HOU_IND", 44.5, -1.5),
    @(13, "SF_CLE",  42.5, -6),
    @(13, "ARI_TB",  48.5,  3.5),
    @(13, "NO_MIA",  44.5,  6),
    @(13, "LA_CAR",  45.5, -3.5),
    @(13, "ATL_NYJ", 43.5, -1.5),
    @(13, "JAX_TEN", 45.5, -1.5),
    @(13, "MIN_SEA", 43.5,  1.5),
    @(13, "LV_LAC",  44.5,  6.5),
    @(13, "BUF_PIT", 46.5, -4.5),
    @(13, "DEN_WAS", 46.5,  2.5),
    @(13, "NYG_NE",  42.5,  3)
)

$startRow = 164
for ($i = 0; $i -lt $week13.Count; $i++) {
    $row = $startRow + $i
    $data = $week13[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
}

$excel.ActiveWindow.ScrollRow = 158
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E165").Select()
